$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Select()

for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = [string]$cell.Value2
    $clean = $text -replace '[\s\u00A0]', '' -replace 'руб\.', ''
    $cell.Value = [double]$clean
}

# column B no longer needs the extra room once the text "руб." is gone,
# so it settles back down to (approximately) its default width
$ws.Columns.Item(2).ColumnWidth = 8.3
